$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Final roster table (rows 2-19), reflecting the new player list / positions / teams.
$data = @(
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Ty Jerome", "PG,SG", "Cleveland Cavaliers"),
    @("Quentin Grimes", "SG,SF", "Dallas Mavericks"),
    @("Kentavious Caldwell-Pope", "SG,SF", "Orlando Magic"),
    @("Harrison Barnes", "SF,PF", "San Antonio Spurs"),
    @("Brandon Miller", "SG,SF", "Charlotte Hornets"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Nick Richards", "C", "Charlotte Hornets"),
    @("Julius Randle", "PF", "Minnesota Timberwolves"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Malcolm Brogdon", "PG,SG", "Washington Wizards"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
